# Apply updated Leve-profit figures (currentAveragePrice / LevePrice / LeveProfit columns)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6324.4546
$ws.Range("J18").Value = 6618.25
$ws.Range("L18").Value = 6618.25
$ws.Range("N18").Value = -7186.25

$ws.Range("H29").Value = 908.3333
$ws.Range("I29").Value = 950
$ws.Range("K29").Value = 2850
$ws.Range("M29").Value = -2569

$ws.Range("H31").Value = 2950
$ws.Range("I31").Value = 2950
$ws.Range("K31").Value = 8850
$ws.Range("M31").Value = -8620

$ws.Range("H33").Value = 1104
$ws.Range("I33").Value = 464.7143
$ws.Range("J33").Value = 1999
$ws.Range("K33").Value = 464.7143
$ws.Range("L33").Value = 1999
$ws.Range("M33").Value = -235.7143
$ws.Range("N33").Value = -2457

$ws.Range("H88").Value = 3269.8
$ws.Range("I88").Value = 3312.25
$ws.Range("J88").Value = 3241.5
$ws.Range("K88").Value = 3312.25
$ws.Range("L88").Value = 3241.5
$ws.Range("M88").Value = -2906.25
$ws.Range("N88").Value = -4053.5

$ws.Range("H91").Value = 3269.8
$ws.Range("I91").Value = 3312.25
$ws.Range("J91").Value = 3241.5
$ws.Range("K91").Value = 3312.25
$ws.Range("L91").Value = 3241.5
$ws.Range("M91").Value = -1908.25
$ws.Range("N91").Value = -6049.5

$ws.Range("H98").Value = 26462.889
$ws.Range("I98").Value = 35431.625
$ws.Range("J98").Value = 13417.454
$ws.Range("K98").Value = 35431.625
$ws.Range("L98").Value = 13417.454
$ws.Range("M98").Value = -33933.625
$ws.Range("N98").Value = -16413.454

$ws.Range("H122").Value = 26462.889
$ws.Range("I122").Value = 35431.625
$ws.Range("J122").Value = 13417.454
$ws.Range("K122").Value = 106294.875
$ws.Range("L122").Value = 40252.362
$ws.Range("M122").Value = -103844.875
$ws.Range("N122").Value = -45152.362

$ws.Range("H125").Value = 3715.24
$ws.Range("I125").Value = 4127.4614
$ws.Range("J125").Value = 3268.6667
$ws.Range("K125").Value = 37147.1526
$ws.Range("L125").Value = 29418.0003
$ws.Range("M125").Value = -34687.1526
$ws.Range("N125").Value = -34338.0003

$ws.Range("H127").Value = 1166.6666
$ws.Range("I127").Value = 1166.6666
$ws.Range("K127").Value = 3499.9998
$ws.Range("M127").Value = 1460.0002

$ws.Range("H131").Value = 4622.1665
$ws.Range("I131").Value = 1933.25
$ws.Range("K131").Value = 5799.75
$ws.Range("M131").Value = -759.75

$ws.Range("H132").Value = 3570.6611
$ws.Range("I132").Value = 3517.4424
$ws.Range("K132").Value = 10552.3272
$ws.Range("M132").Value = -8022.3272

$ws.Range("H138").Value = 3718.0952
$ws.Range("I138").Value = 2735.182
$ws.Range("K138").Value = 8205.545999999998
$ws.Range("M138").Value = -3065.545999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4368.114
$ws.Range("I32").Value = 4237.1396
$ws.Range("K32").Value = 4237.1396
$ws.Range("M32").Value = -3950.1396

$ws.Range("H44").Value = 27245
$ws.Range("J44").Value = 27245
$ws.Range("L44").Value = 27245
$ws.Range("N44").Value = -28221

$ws.Range("H111").Value = 92500
$ws.Range("J111").Value = 92500
$ws.Range("L111").Value = 92500
$ws.Range("N111").Value = -100680

$ws.Range("H130").Value = 42166.332
$ws.Range("J130").Value = 42166.332
$ws.Range("L130").Value = 42166.332
$ws.Range("N130").Value = -52206.332

$ws.Range("H132").Value = 7234.5293
$ws.Range("I132").Value = 8123.773
$ws.Range("K132").Value = 24371.319
$ws.Range("M132").Value = -21841.319

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 49755.25
$ws.Range("I26").Value = 50720.285
$ws.Range("J26").Value = 43000
$ws.Range("K26").Value = 50720.285
$ws.Range("L26").Value = 43000
$ws.Range("M26").Value = -50428.285
$ws.Range("N26").Value = -43584

$ws.Range("H81").Value = 19780
$ws.Range("J81").Value = 19780
$ws.Range("L81").Value = 19780
$ws.Range("N81").Value = -21902

$ws.Range("H84").Value = 19780
$ws.Range("J84").Value = 19780
$ws.Range("L84").Value = 59340
$ws.Range("N84").Value = -69948

$ws.Range("H96").Value = 12500
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = $null

$ws.Range("H105").Value = 105049.9
$ws.Range("I105").Value = 145214.14
$ws.Range("K105").Value = 145214.14
$ws.Range("M105").Value = -143467.14

$ws.Range("H134").Value = 10822.071
$ws.Range("I134").Value = 11628.72
$ws.Range("K134").Value = 34886.16
$ws.Range("M134").Value = -32351.16

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1111.6471
$ws.Range("I22").Value = 604.375
$ws.Range("K22").Value = 604.375
$ws.Range("M22").Value = -254.375

$ws.Range("H86").Value = 14603.25
$ws.Range("J86").Value = 15278.1
$ws.Range("L86").Value = 15278.1
$ws.Range("N86").Value = -17524.1

$ws.Range("H89").Value = 14603.25
$ws.Range("J89").Value = 15278.1
$ws.Range("L89").Value = 76390.5
$ws.Range("N89").Value = -87622.5

$ws.Range("H99").Value = 5534780.5
$ws.Range("I99").Value = 12905293
$ws.Range("K99").Value = 12905293
$ws.Range("M99").Value = -12903795

$ws.Range("H114").Value = 3000
$ws.Range("J114").Value = 3000
$ws.Range("L114").Value = 3000
$ws.Range("N114").Value = -11678

$ws.Range("H126").Value = 5534780.5
$ws.Range("I126").Value = 12905293
$ws.Range("K126").Value = 38715879
$ws.Range("M126").Value = -38713409

$ws.Range("H132").Value = 25636.176
$ws.Range("I132").Value = 1854.3334
$ws.Range("J132").Value = 204000
$ws.Range("K132").Value = 5563.0002
$ws.Range("L132").Value = 612000
$ws.Range("M132").Value = -3033.0002
$ws.Range("N132").Value = -617060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1634869.6
$ws.Range("I121").Value = 1808879.2
$ws.Range("J121").Value = 1541172.1
$ws.Range("K121").Value = 5426637.6
$ws.Range("L121").Value = 4623516.300000001
$ws.Range("M121").Value = -5425327.6
$ws.Range("N121").Value = -4626136.300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 53110.1
$ws.Range("J105").Value = 53110.1
$ws.Range("L105").Value = 53110.1
$ws.Range("N105").Value = -60098.1

$ws.Range("H122").Value = 5419.7607
$ws.Range("I122").Value = 3478.1943
$ws.Range("K122").Value = 10434.5829
$ws.Range("M122").Value = -7984.582900000001

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("N132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16830.916
$ws.Range("I7").Value = 26007.055
$ws.Range("J7").Value = 7654.778
$ws.Range("K7").Value = 26007.055
$ws.Range("L7").Value = 7654.778
$ws.Range("M7").Value = -25895.055
$ws.Range("N7").Value = -7878.778

$ws.Range("H16").Value = 5182.2666
$ws.Range("I16").Value = 5574.7
$ws.Range("J16").Value = 4397.4
$ws.Range("K16").Value = 5574.7
$ws.Range("L16").Value = 4397.4
$ws.Range("M16").Value = -5404.7
$ws.Range("N16").Value = -4737.4

$ws.Range("H40").Value = 28165.334
$ws.Range("I40").Value = 47444.668
$ws.Range("J40").Value = 13705.833
$ws.Range("K40").Value = 47444.668
$ws.Range("L40").Value = 13705.833
$ws.Range("M40").Value = -47308.668
$ws.Range("N40").Value = -13977.833

$ws.Range("H122").Value = 4865.9165
$ws.Range("I122").Value = 4352.3335
$ws.Range("K122").Value = 13057.0005
$ws.Range("M122").Value = -10607.0005

$ws.Range("H126").Value = 16830.916
$ws.Range("I126").Value = 26007.055
$ws.Range("J126").Value = 7654.778
$ws.Range("K126").Value = 78021.16500000001
$ws.Range("L126").Value = 22964.334
$ws.Range("M126").Value = -75551.16500000001
$ws.Range("N126").Value = -27904.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 14398.25
$ws.Range("I81").Value = 20639.4
$ws.Range("K81").Value = 41278.8
$ws.Range("M81").Value = -40217.8

$ws.Range("H84").Value = 14398.25
$ws.Range("I84").Value = 20639.4
$ws.Range("K84").Value = 206394
$ws.Range("M84").Value = -201090

$ws.Range("H96").Value = 3433
$ws.Range("I96").Value = 2667.6667
$ws.Range("J96").Value = 4198.3335
$ws.Range("K96").Value = 2667.6667
$ws.Range("L96").Value = 4198.3335
$ws.Range("M96").Value = -1294.6667
$ws.Range("N96").Value = -6944.3335

$ws.Range("H99").Value = 41716
$ws.Range("I99").Value = 23432
$ws.Range("K99").Value = 23432
$ws.Range("M99").Value = -20437

$ws.Range("H100").Value = 72685.625
$ws.Range("I100").Value = 43580.832
$ws.Range("K100").Value = 87161.664
$ws.Range("M100").Value = -86620.664

$ws.Range("H113").Value = 1210.7106
$ws.Range("J113").Value = 2431.0833
$ws.Range("L113").Value = 7293.249899999999
$ws.Range("N113").Value = -11633.2499

$ws.Range("H122").Value = 5218.7144
$ws.Range("I122").Value = 2411.842
$ws.Range("J122").Value = 8551.875
$ws.Range("K122").Value = 7235.526
$ws.Range("L122").Value = 25655.625
$ws.Range("M122").Value = -4785.526
$ws.Range("N122").Value = -30555.625

$ws.Range("H126").Value = 20558.348
$ws.Range("I126").Value = 27096.562
$ws.Range("K126").Value = 81289.686
$ws.Range("M126").Value = -78819.686
